$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add P1=14, Q1=15, matching style of existing header cells ---
$ws.Cells.Item(1, 16).Value = 14   # P1
$ws.Cells.Item(1, 17).Value = 15   # Q1

foreach ($col in @(16, 17)) {
    $c = $ws.Cells.Item(1, $col)
    $c.Font.Bold = $true
    $c.Borders.LineStyle = 1
    $c.HorizontalAlignment = -4108   # xlCenter
    $c.VerticalAlignment = -4160     # xlTop
}

# --- Data rows 2..25 ---
# Columns I, K, M, O get their values swapped (1<->2), and new columns P, Q = 2 are appended.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P (new)
    $ws.Cells.Item($r, 17).Value = 2   # Q (new)
}
